# aggiornamento fino a 9 agosto 2021
# Append new daily rows (329-343, dates 2021-07-26 .. 2021-08-09) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries a date-formatted style (same as the existing A-column cells).
# Copy that formatting down onto the new A-cells before writing their values.
$ws.Range("A328").Copy($ws.Range("A329:A343"))

$rows = @(
    @(329, 44403, 1, 2, 28.44950213371266),
    @(330, 44404, 0, 2, 28.44950213371266),
    @(331, 44405, 0, 2, 28.44950213371266),
    @(332, 44406, 0, 2, 28.44950213371266),
    @(333, 44407, 0, 2, 28.44950213371266),
    @(334, 44408, 0, 2, 28.44950213371266),
    @(335, 44409, 0, 1, 14.22475106685633),
    @(336, 44410, 0, 0, 0),
    @(337, 44411, 0, 0, 0),
    @(338, 44412, 0, 0, 0),
    @(339, 44413, 0, 0, 0),
    @(340, 44414, 1, 1, 14.22475106685633),
    @(341, 44415, 0, 1, 14.22475106685633),
    @(342, 44416, 2, 3, 42.67425320056899),
    @(343, 44417, 1, 4, 56.89900426742533)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
